$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace speaker label "G4" with "T" in the Speaker column (D)
$speakerRows = @(2, 4, 13, 19, 22, 23, 27, 33, 38, 40, 44, 46, 48, 55, 57, 59, 60)
foreach ($r in $speakerRows) {
    $ws.Range("D$r").Value = "T"
}

# Replace inline mentions of "G4" with "T" inside Sentence column (E)
$ws.Range("E6").Value = "[pointing his finger to T] I got that on lock down."
$ws.Range("E20").Value = "[turning his head towards T] Eighty-four?"
$ws.Range("E36").Value = "[pointing something in his paper to T who is coming around to his  side] I know that there is some pattern… between this one and this one."
$ws.Range("E58").Value = "[to T] Have you guys got more problems?"
